$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the 11 "same character" diagonal cells: distance 0 -> 0.5, ---
# --- and give them the highlighted/bordered look (already-used style). ---
$diag = @("C3", "D4", "E5", "F6", "G7", "H8", "I9", "J10", "K11", "L12", "M13")
foreach ($addr in $diag) {
    $rng = $ws.Range($addr)
    $rng.Value = 0.5
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
}

# --- Add a small legend row describing that "same character" distance is 0.5 ---
$ws.Range("O5:Q5").Merge()
$ws.Range("O5").Value = "같은 글자에 대한 거리"
$ws.Range("O5:Q5").Borders.LineStyle = 1
$ws.Range("O5:Q5").HorizontalAlignment = -4108
$ws.Range("R5").Value = 0.5
$ws.Range("R5").Borders.LineStyle = 1
$ws.Range("R5").HorizontalAlignment = -4108

# --- Leave the cursor where the author left it when they saved ---
$ws.Range("V15").Select() | Out-Null
